$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the row above onto the new row, then fill values.
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial()

$ws.Range("A5").Value = 42608.901620370372
$ws.Range("B5").Value = 77
$ws.Range("C5:M5").Value = 0
$ws.Range("N5").Value = "Random"
